# add multiline support for description box
#
# The change splits a handful of single runs of text into multiple runs
# that carry the *same* run formatting (rPr) as the original run. A plain
# Range.Text / Find-Replace write collapses back into one run when the
# neighbouring run(s) end up with identical formatting, so instead we
# pin each future run boundary with a transient bookmark *before* we
# mutate the text there. Bookmarks are zero-width marks that keep the
# text on either side of them from being re-coalesced into a single run
# by the engine's run-merge normalisation, so once the edit is done we
# can just delete the scratch bookmarks (the split already happened) -
# except for the one spot that keeps Word's own "_GoBack" bookmark.

function Split-FieldWord($d, $paraIndex, $newParts, $useGoBack) {
    # $newParts describes how the *first character* of the paragraph's
    # (single-run) text is rewritten, followed by the untouched remainder:
    #   2 parts -> [replacement-for-old-char-1, untouched-rest]
    #   3 parts -> [replacement-a, replacement-b, untouched-rest]
    #              (old char 1 becomes replacement-a + replacement-b,
    #               each landing in its own run)
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $start = $r.Start

    # Mark the boundary between old char 1 and the untouched remainder
    # before editing, so that boundary survives as a run split.
    $bm1Name = "ZZZ_TMP_SPLIT_1"
    $d.Bookmarks.Add($bm1Name, $d.Range($start + 1, $start + 1)) | Out-Null

    if ($newParts.Length -eq 2) {
        $rngC = $d.Range($start, $start + 1)
        $rngC.Text = $newParts[0]
        $d.Bookmarks($bm1Name).Delete()
    }
    elseif ($newParts.Length -eq 3) {
        $combined = $newParts[0] + $newParts[1]
        $rngC = $d.Range($start, $start + 1)
        $rngC.Text = $combined
        $d.Bookmarks($bm1Name).Delete()

        # Mark the boundary between replacement-a and replacement-b,
        # added only now (after the text landed), so it sits between
        # the two new characters rather than around them.
        $bm2Pos = $start + $newParts[0].Length
        if ($useGoBack) {
            $d.Bookmarks.Add("_GoBack", $d.Range($bm2Pos, $bm2Pos)) | Out-Null
        }
        else {
            $d.Bookmarks.Add("ZZZ_TMP_SPLIT_2", $d.Range($bm2Pos, $bm2Pos)) | Out-Null
            $d.Bookmarks("ZZZ_TMP_SPLIT_2").Delete()
        }
    }
    else {
        throw "Split-FieldWord: unsupported newParts length $($newParts.Length)"
    }
}

$d = $word.ActiveDocument

# --- 1. "Login_user(login,password) " -> "l" + "ogin_user(login,password) "
$findRng = $d.Content
$findRng.Find.Execute("Login_user(login,password) ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$luStart = $findRng.Start
$d.Bookmarks.Add("ZZZ_TMP_SPLIT_1", $d.Range($luStart + 1, $luStart + 1)) | Out-Null
$luRng = $d.Range($luStart, $luStart + 1)
$luRng.Text = "l"
$d.Bookmarks("ZZZ_TMP_SPLIT_1").Delete()

# --- 2. database field-name bullet list: prefix each field name with "_"
#         and lower-case its first letter, e.g. "Id" -> "_id", "Title" -> "_title"
Split-FieldWord $d 50 @("_i", "d") $false
Split-FieldWord $d 51 @("_", "t", "itle") $false
Split-FieldWord $d 52 @("_", "d", "irector") $false
Split-FieldWord $d 53 @("_", "d", "escription") $false
Split-FieldWord $d 54 @("_", "d", "uration") $false
Split-FieldWord $d 55 @("_", "a", "ctors") $false
Split-FieldWord $d 56 @("_", "g", "enre") $false
Split-FieldWord $d 57 @("_", "i", "con_path") $false
Split-FieldWord $d 58 @("_", "a", "vg_rate") $false
Split-FieldWord $d 59 @("_", "n", "umber_of_users") $false

Split-FieldWord $d 79 @("_", "u", "ser_id") $false
Split-FieldWord $d 80 @("_", "m", "ovie_id") $false
Split-FieldWord $d 81 @("_", "s", "core") $false

Split-FieldWord $d 88 @("_", "l", "ogin") $false
Split-FieldWord $d 89 @("_", "e", "mail") $false
Split-FieldWord $d 90 @("_", "p", "assword") $false
Split-FieldWord $d 91 @("_", "i", "con_path") $false

# --- 3. move the document's "_GoBack" bookmark from the end of the test
#         summary paragraph to the middle of the last "Id" -> "_id" edit.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

Split-FieldWord $d 92 @("_", "i", "d") $true

Write-Host "done"
